# This script updates the "symbol list" crypto price table on Sheet1.
# It refreshes the Price (col D) and Volume(1h) (col E) figures for most
# rows, and shifts rows 9-14 up by one slot (a new coin "One" was
# inserted at the top of that block, pushing WazirX, MandalaExchangeToken,
# BitrueCoin, BitMartToken and BitForexToken down one row each), updating
# their Coin (col B) / Link (col C) / Price / Volume values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, taken from the published diff.
$updates = [ordered]@{
    'D2'  = '258.27'
    'E2'  = '4.98%'
    'D3'  = '27.25'
    'E3'  = '-3.70%'
    'D4'  = '5.214'
    'E4'  = '-1.50%'
    'D5'  = '0.05922'
    'E5'  = '3.69%'
    'D6'  = '6.705'
    'E6'  = '0.94%'
    'D7'  = '0.8668'
    'E7'  = '0.35%'
    'D8'  = '1.003'
    'E8'  = '13.35%'
    'B9'  = 'One'
    'C9'  = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D9'  = '0.0006074'
    'E9'  = '-94.03%'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1417'
    'E10' = '2.10%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.07184'
    'E11' = '1.29%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.03147'
    'E12' = '-0.16%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09221'
    'E13' = '-0.14%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001557'
    'E14' = '1.98%'
    'D15' = '0.005903'
    'E15' = '-1.99%'
    'D16' = '3.499'
    'E16' = '0.06%'
    'D17' = '3.268'
    'E17' = '1.64%'
    'D18' = '2.226'
    'E18' = '2.48%'
    'E19' = '-0.72%'
    'D20' = '0.03561'
    'E20' = '6.31%'
    'E21' = '-0.26%'
    'D22' = '3.522'
    'E22' = '0.98%'
    'D23' = '0.04187'
    'E23' = '2.10%'
    'E24' = '1.43%'
    'D25' = '0.001218'
    'E25' = '-0.02%'
    'E26' = '8.70%'
    'E27' = '-0.03%'
    'E28' = '2.64%'
    'D40' = '0.03835'
    'E40' = '1.11%'
    'D41' = '0.006573'
    'E41' = '75.27%'
    'D42' = '0.1105'
    'E42' = '3.30%'
    'D43' = '0.002198'
    'E43' = '-10.23%'
    'D44' = '0.01084'
    'E44' = '14.59%'
    'D45' = '0.00005400'
    'E45' = '2.50%'
    'E46' = '-0.08%'
    'E47' = '22.34%'
    'E48' = '-1.29%'
    'E49' = '-0.08%'
    'E50' = '-0.08%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "258.27") and
    # percent-looking strings (e.g. "4.98%") are preserved verbatim as
    # text instead of being auto-converted to numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
